$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2 through 40 (columns A = consumption, B = timestamp serial)
$data = @(
    @(5376, 45797),
    @(5355, 45797.01041666666),
    @(5348, 45797.02083333334),
    @(5284, 45797.03125),
    @(5230, 45797.04166666666),
    @(5221, 45797.05208333334),
    @(5166, 45797.0625),
    @(5197, 45797.07291666666),
    @(5159, 45797.08333333334),
    @(5131, 45797.09375),
    @(5124, 45797.10416666666),
    @(5163, 45797.11458333334),
    @(5247, 45797.125),
    @(5190, 45797.13541666666),
    @(5249, 45797.14583333334),
    @(5234, 45797.15625),
    @(5337, 45797.16666666666),
    @(5309, 45797.17708333334),
    @(5342, 45797.1875),
    @(5347, 45797.19791666666),
    @(5505, 45797.20833333334),
    @(5582, 45797.21875),
    @(5686, 45797.22916666666),
    @(5735, 45797.23958333334),
    @(5920, 45797.25),
    @(5975, 45797.26041666666),
    @(6003, 45797.27083333334),
    @(5979, 45797.28125),
    @(6017, 45797.29166666666),
    @(6000, 45797.30208333334),
    @(5899, 45797.3125),
    @(5807, 45797.32291666666),
    @(5652, 45797.33333333334),
    @(5556, 45797.34375),
    @(5471, 45797.35416666666),
    @(5371, 45797.36458333334),
    @(5346, 45797.375),
    @(5320, 45797.38541666666),
    @(5220, 45797.39583333334)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

# Rows 27-40 are brand new; make sure column B keeps the same timestamp
# number format as the existing rows (copy format from B2).
$ws.Range("B27:B40").NumberFormat = $ws.Range("B2").NumberFormat
